$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Drop the stray "_GoBack" bookmark that used to sit between "BLANC, " and
#    "LEMPEREUR" in the authors' byline. Word regenerates bookmark ids from
#    their position in the saved XML, so removing this one automatically
#    shifts the five "_Toc..." bookmark ids (1-5) down to (0-4).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. The conclusion paragraph loses its trailing clause "pour les envoyer en
#    BDD", and Word's cursor-position bookmark ("_GoBack") ends up right
#    after "capteurs " instead. Locate that spot *before* touching the
#    trailing text (so the insertion point is unambiguous), drop the new
#    bookmark there, and only then delete the trailing words. Doing it in
#    this order keeps the new bookmark outside of the range being edited so
#    it survives the edit instead of being swallowed by it.
# ---------------------------------------------------------------------------
$marker = "En conclusion nous avons réussi à nous connecté et à recouper les valeurs des capteurs "
$hit = $d.Content
$hit.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertionPoint = $d.Range($hit.End, $hit.End)
$d.Bookmarks.Add("_GoBack", $insertionPoint)

$tail = "pour les envoyer en BDD"
$d.Content.Find.Execute($tail, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
